$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $text = $val.ToString()
        if ($text.Contains(",")) {
            $parts = $text.Split(",")
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }
            $n = $trimmed.Length
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $trimmed[$i]
            }
            $newVal = [string]::Join(", ", $reversed)
            $cell.Value2 = $newVal
        }
    }
}
